$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CompanyName","URL1","URL2","URL3","URL4"),
    @("Dell","https://www.dell.com/en-in/shop/deals/laptop-deals","https://www.dell.com/support/home/en-in","https://www.dell.com/support/home/en-in?app=drivers","https://www.dell.com/en-in/shop/laptops-2-in-1-pcs/sc/laptops"),
    @("Techmahindra","https://careers.techmahindra.com/","https://careers.techmahindra.com/CurrentOpportunity.aspx","https://www.techmahindra.com/en-in/search/?key=jobs","https://www.techmahindra.com/en-in/techmahindra-overview/"),
    @("Verizon","https://secure.verizon.com/vzauth/UI/Login","https://www.verizon.com/smartphones/","https://www.verizon.com/about/careers","https://www.verizon.com/plans/"),
    @("infosys","https://www.infosys.com/careers/apply.html","https://career.infosys.com/","https://www.infosys.com/about/history.html","https://www.infosys.com/about.html"),
    @("Amazon","https://www.amazon.in/Home-Kitchen/b?ie=UTF8&node=976442031","https://www.amazon.in/electronics/b?ie=UTF8&node=976419031","https://www.amazon.in/b?ie=UTF8&node=7459780031","https://www.amazon.in/amazonpay/home?ref_=apay_logo_APayDashboard"),
    @("Wipro","https://careers.wipro.com/careers-home/","https://www.wipro.com/location/","https://www.wipro.com/about-us/","https://www.wipro.com/en-US/"),
    @("Google","https://translate.google.co.in/","https://www.google.co.in/maps","https://books.google.co.in/","https://news.google.co.in/"),
    @("Microsoft","https://account.microsoft.com/account/Account","https://www.microsoft.com/en-in/download","https://careers.microsoft.com/","https://apps.microsoft.com/")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws.Range("C11").Select()

$wb.Save()
